# Updated capital structure database
# Applies refreshed financial metrics to the Finland / Drugs (Biotechnology)
# rows (2-4) in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
# ebitdard_margin / ebitda_margin / operating_margin / after_tax_operating_margin
# and net_margin no longer apply to this row; drop those cells entirely.
$ws.Range("G2:J2").ClearContents()
$ws.Range("L2").ClearContents()

$ws.Range("K2").Value = -18.47

$ws.Range("U2").Value  = 18.75
$ws.Range("V2").Value  = 0.091552734375
$ws.Range("W2").Value  = 3.695061728395062
$ws.Range("X2").Value  = 0.06266027892055928
$ws.Range("Y2").Value  = 3.632401449474502
$ws.Range("Z2").Value  = 0
$ws.Range("AA2").Value = 3.137623545295393
$ws.Range("AB2").Value = 0.0593555342464834
$ws.Range("AC2").Value = 3.07826801104891
$ws.Range("AD2").Value = 9.949999999999999
$ws.Range("AE2").Value = 0.2076740687619426
$ws.Range("AF2").Value = 10.15767406876194
$ws.Range("AG2").Value = -8.592325931238058
$ws.Range("AH2").Value = 0.04725429837649175
$ws.Range("AI2").Value = 0.3801855671500312
$ws.Range("AJ2").Value = -0.04379199729072185
$ws.Range("AK2").Value = -1.078398269944942
$ws.Range("AL2").Value = 0.46
$ws.Range("AM2").Value = 0.4309999999999999
$ws.Range("AN2").Value = -0.5674689175316527
$ws.Range("AO2").Value = -39.06521739130435
$ws.Range("AP2").Value = 0.4900379794250062
$ws.Range("AQ2").Value = -41.69373549883991

# --- Row 3 ---------------------------------------------------------------
$ws.Range("G3:J3").ClearContents()
$ws.Range("L3").ClearContents()

$ws.Range("K3").Value = -16

$ws.Range("U3").Value  = 13.1
$ws.Range("V3").Value  = 0.07867867867867867
$ws.Range("W3").Value  = 8
$ws.Range("X3").Value  = 0.05952898253473292
$ws.Range("Y3").Value  = 7.940471017465267
$ws.Range("Z3").Value  = 0
$ws.Range("AA3").Value = 6.506444046985285
$ws.Range("AB3").Value = 0.05882196402581642
$ws.Range("AC3").Value = 6.447622082959469
$ws.Range("AD3").Value = 3.12
$ws.Range("AE3").Value = 0.2076740687619426
$ws.Range("AF3").Value = 3.327674068761943
$ws.Range("AG3").Value = -9.772325931238058
$ws.Range("AH3").Value = 0.01959441585129755
$ws.Range("AI3").Value = 0.2879190093927266
$ws.Range("AJ3").Value = -0.06235226796609382
$ws.Range("AK3").Value = 6.336096497705649
$ws.Range("AL3").Value = 0.352
$ws.Range("AM3").Value = 0.323
$ws.Range("AN3").Value = -0.2036021926389976
$ws.Range("AO3").Value = -44.31818181818182
$ws.Range("AP3").Value = 0.6377137778150651
$ws.Range("AQ3").Value = -48.29721362229103

# --- Row 4 ---------------------------------------------------------------
# (no ebitdard/ebitda/operating-margin or net_margin cells in this row to begin with)
$ws.Range("K4").Value = -2.47

$ws.Range("U4").Value  = 5.65
$ws.Range("V4").Value  = 0.1475195822454308
$ws.Range("W4").Value  = -0.6098765432098766
$ws.Range("X4").Value  = 0.06579157530638566
$ws.Range("Y4").Value  = -0.6756681185162623
$ws.Range("AA4").Value = -0.2311969563944981
$ws.Range("AB4").Value = 0.0598891044671504
$ws.Range("AC4").Value = -0.2910860608616485
$ws.Range("AD4").Value = 6.83
$ws.Range("AF4").Value = 6.83
$ws.Range("AG4").Value = 1.18
$ws.Range("AH4").Value = 0.1513405716818081
$ws.Range("AI4").Value = 0.4505277044854881
$ws.Range("AJ4").Value = 0.02988855116514691
$ws.Range("AK4").Value = 0.1240799158780231
$ws.Range("AL4").Value = 0.108
$ws.Range("AM4").Value = 0.108
$ws.Range("AN4").Value = -3.090497737556561
$ws.Range("AO4").Value = -21.94444444444445
$ws.Range("AP4").Value = -0.5339366515837103
$ws.Range("AQ4").Value = -21.94444444444445
